# This script normalizes vaccine/product name labels in the shared strings used
# across all four worksheets:
#   - Strip trailing footnote markers like " [1]", " [2]", " [3]", " [4]", " [5]"
#     (the preceding space is kept).
#   - Collapse embedded line breaks within a label into a single space so that
#     multi-line labels (e.g. "Fluzone\nQuadrivalent") become one line
#     (e.g. "Fluzone Quadrivalent").
# Setting the cleaned value on every affected cell causes Excel to rebuild the
# shared-string table accordingly (including merging the "Afluria\nQuadrivalent"
# duplicate into the already-existing single-line "Afluria Quadrivalent" string).

$wb = $excel.ActiveWorkbook

# Sheet: Pediatric VFC Vaccine 
$ws = $wb.Worksheets.Item('Pediatric VFC Vaccine ')
$updates = @{
    'A2' = 'DTaP '
    'A3' = 'DTaP '
    'A4' = 'DTaP '
    'A5' = 'DTaP-IPV '
    'A6' = 'DTaP-IPV '
    'A7' = 'DTaP-IPV '
    'A8' = 'DTaP-Hep B-IPV '
    'A9' = 'DTaP-IP-HI '
    'A10' = 'e-IPV '
    'A11' = 'Hepatitis A Pediatric '
    'A12' = 'Hepatitis A Pediatric '
    'A13' = 'Hepatitis A-Hepatitis B 18 only '
    'A14' = 'Hepatitis B  Pediatric/Adolescent'
    'A15' = 'Hepatitis B  Pediatric/Adolescent'
    'B15' = 'Recombivax HB'
    'A16' = 'Hepatitis B  Pediatric/Adolescent'
    'B16' = 'Recombivax HB'
    'A17' = 'Hib '
    'A18' = 'Hib '
    'A19' = 'Hib '
    'A20' = 'HPV - Human Papillomavirus 9-valent '
    'A21' = 'MENB - Meningococcal Group B '
    'A22' = 'MENB - Meningococcal Group B '
    'A23' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A24' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A25' = 'Measles, Mumps and Rubella (MMR) '
    'A26' = 'MMR/Varicella '
    'A27' = 'Pneumococcal 13-valent  (Pediatric)'
    'A29' = 'Rotavirus, Live, Oral, Pentavalent '
    'A30' = 'Rotavirus, Live, Oral, Pentavalent '
    'A31' = 'Rotavirus, Live, Oral, Oral '
    'A32' = 'Tetanus and Diphtheria Toxoids '
    'A33' = 'Tetanus and Diphtheria Toxoids '
    'A34' = 'Tetanus and Diphtheria Toxoids '
    'A35' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A36' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A37' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A38' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A39' = 'Varicella '
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Sheet: Adult Vaccine 
$ws = $wb.Worksheets.Item('Adult Vaccine ')
$updates = @{
    'A2' = 'Hepatitis A Adult '
    'A3' = 'Hepatitis A Adult '
    'A4' = 'Hepatitis A Adult '
    'A5' = 'Hepatitis A-Hepatitis B Adult '
    'A6' = 'Hepatitis B Adult '
    'A7' = 'Hepatitis B Adult '
    'A8' = 'Hepatitis B Adult '
    'A9' = 'HPV-Human Papillomavirus 9 Valent '
    'A10' = 'Measles, Mumps,  Rubella '
    'A11' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A12' = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
    'A13' = 'MENB - Meningococcal Group B '
    'A14' = 'MENB - Meningococcal Group B '
    'A15' = 'MENB - Meningococcal Group B '
    'A16' = 'Pneumococcal 13-valent '
    'A19' = 'Tetanus and Diphtheria Toxoids '
    'A20' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A21' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A22' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A23' = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
    'A24' = 'Varicella '
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Sheet: Pediatric Influenza Vaccine 
$ws = $wb.Worksheets.Item('Pediatric Influenza Vaccine ')
$updates = @{
    'A2' = 'Influenza  (Age 6 months and older)'
    'B2' = 'Fluzone Quadrivalent'
    'A3' = 'Influenza  (Age 6-35 months)'
    'B3' = 'Fluzone Quadrivalent Pediatric dose'
    'A4' = 'Influenza  (Age 6 months and older)'
    'B4' = 'Fluzone Quadrivalent'
    'A5' = 'Influenza  (Age 6 months and older)'
    'B5' = 'Fluzone Quadrivalent'
    'A6' = 'Influenza  (Age 6 months and older)'
    'B6' = 'Fluarix Quadrivalent'
    'A7' = 'Influenza  (Age 6 months and older)'
    'B7' = 'FluLaval Quadrivalent'
    'A8' = 'Influenza  (Age 6 months and older)'
    'B8' = 'FluLaval Quadrivalent'
    'A9' = 'Influenza  (Age 4 years and older)'
    'A10' = 'Influenza  (Age 4 years and older)'
    'A11' = 'Influenza  (Age 6 -35 months)'
    'A12' = 'Influenza  (Age 36 months and older)'
    'A13' = 'Influenza  (Age 6 months and older)'
    'A14' = 'Influenza  Live, Intranasal (Age 2-49 years)'
    'B14' = 'FluMist Quadrivalent'
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Sheet: Adult Influenza Vaccine 
$ws = $wb.Worksheets.Item('Adult Influenza Vaccine ')
$updates = @{
    'A2' = 'Influenza  (Age 6 months and older)'
    'B2' = 'Fluzone Quadrivalent'
    'A3' = 'Influenza  (Age 6 months and older)'
    'B3' = 'Fluzone Quadrivalent'
    'A4' = 'Influenza  (Age 6 months and older)'
    'B4' = 'Fluzone Quadrivalent'
    'A5' = 'Influenza  (Age 6 months and older)'
    'B5' = 'Fluarix Quadrivalent'
    'A6' = 'Influenza  (Age 6 months and older)'
    'B6' = 'FluLaval Quadrivalent'
    'A7' = 'Influenza  (Age 6 months and older)'
    'B7' = 'FluLaval Quadrivalent'
    'A8' = 'Influenza  (Age 4 years and older)'
    'A9' = 'Influenza  (Age 4 years and older)'
    'A10' = 'Influenza  (Age 36 months and older)'
    'B10' = 'Afluria Quadrivalent'
    'A11' = 'Influenza  (Age 6 months and older)'
    'B11' = 'Afluria Quadrivalent'
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
